$wb = $excel.ActiveWorkbook
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before the "总计" sheet
# ---------------------------------------------------------------
# NOTE: sheet references returned by Worksheets.Item(...) track sheet
# POSITION, not identity. Once Worksheets.Add() inserts a sheet in
# front of "总计", any reference grabbed beforehand ends up pointing at
# the newly inserted sheet instead. So we re-resolve "总计" by name
# *after* the insert, right before we need to edit it (below).
$totalSheetBeforeInsert = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBeforeInsert)
$newSheet.Name = "2022-Q1"

# Copy header-row (B1:H1) + index-column (A2:A5) cell formatting from the
# "2021-Q4" sheet, which uses the same style (s="2") as the "总计" sheet.
$templateSheet.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2:A5").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# Headers
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Columns B-G hold numeric-looking strings that must stay text (leading
# zeros in fund codes, fixed decimal formatting, etc.)
$newSheet.Range("B2:G5").NumberFormat = "@"

# Row 2
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "011501"
$newSheet.Cells.Item(2,3).Value = "方正富邦汇福一年定期开放灵活配置混合A"
$newSheet.Cells.Item(2,4).Value = "4.76"
$newSheet.Cells.Item(2,5).Value = "41.05"
$newSheet.Cells.Item(2,6).Value = "3.02"
$newSheet.Cells.Item(2,7).Value = "0.1438"
$newSheet.Cells.Item(2,8).Value = 4

# Row 3
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "013714"
$newSheet.Cells.Item(3,3).Value = "方正富邦泰利12个月持有期混合A"
$newSheet.Cells.Item(3,4).Value = "3.60"
$newSheet.Cells.Item(3,5).Value = "20.66"
$newSheet.Cells.Item(3,6).Value = "0.85"
$newSheet.Cells.Item(3,7).Value = "0.0306"
$newSheet.Cells.Item(3,8).Value = 5

# Row 4
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "011502"
$newSheet.Cells.Item(4,3).Value = "方正富邦汇福一年定期开放灵活配置混合C"
$newSheet.Cells.Item(4,4).Value = "0.09"
$newSheet.Cells.Item(4,5).Value = "41.05"
$newSheet.Cells.Item(4,6).Value = "3.02"
$newSheet.Cells.Item(4,7).Value = "0.0027"
$newSheet.Cells.Item(4,8).Value = 4

# Row 5
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "013715"
$newSheet.Cells.Item(5,3).Value = "方正富邦泰利12个月持有期混合C"
$newSheet.Cells.Item(5,4).Value = "0.10"
$newSheet.Cells.Item(5,5).Value = "20.66"
$newSheet.Cells.Item(5,6).Value = "0.85"
$newSheet.Cells.Item(5,7).Value = "0.0008"
$newSheet.Cells.Item(5,8).Value = 5

# ---------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: prepend a 2022-Q1 row and
#    shift the existing rows down, renumbering the index column.
# ---------------------------------------------------------------
# Re-resolve by name now that the sheet collection has changed (see note above).
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# A2 should carry the same index-column style as the rest of the rows.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# B2:D2 should be plain (unstyled), matching B3:D3.
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 4
$totalSheet.Cells.Item(2,4).Value = 0.18

$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2

# ---------------------------------------------------------------
# 3. Restore the original active sheet ("2021-Q2") so the edit doesn't
#    leave an unrelated side effect on tab selection.
# ---------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
